$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update swap-to pin population values for Positronic #1 and #9
$ws.Range("D2").Value = 6
$ws.Range("D10").Value = 7

# Add a cautionary note in a new column, referenced from I2
$ws.Range("I2").Value = "Note: the user should also be careful that the swapped Positronics have the SAME PINS (not just total number) populated"

# Give H1 the same header formatting as the rest of row 1 (bold, centered)
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Move the selection to the newly added note cell
$ws.Range("I2").Select()
